$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new values (plain decimals like "567.75") would otherwise be auto-parsed
# as numbers by the Value setter, so those cells are pre-formatted as Text (@)
# to preserve the literal string, matching the source inline-string cells.

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Rows 35/36 swap identity: row 35 becomes FirstDigitalUSD (was ImmutableX),
# row 36 becomes ImmutableX (was FirstDigitalUSD); plus Price / Volume(1h) refresh
# for every row in the table.
$ws.Range("D2").Value = '61.363.58'
$ws.Range("E2").Value = '  +0.29%  '
$ws.Range("D3").Value = '2.407.19'
$ws.Range("E3").Value = '  -0.45%  '
$ws.Range("E4").Value = '  +0.53%  '
$ws.Range("D5").Value = '567.75'
$ws.Range("E5").Value = '  -0.21%  '
$ws.Range("D6").Value = '142.85'
$ws.Range("E6").Value = '  +2.05%  '
$ws.Range("E7").Value = '  -0.32%  '
$ws.Range("D8").Value = '0.535'
$ws.Range("E8").Value = '  -0.36%  '
$ws.Range("D9").Value = '2.417.92'
$ws.Range("E9").Value = '  +0.60%  '
$ws.Range("D10").Value = '0.108'
$ws.Range("E10").Value = '  +1.45%  '
$ws.Range("E11").Value = '  +0.35%  '
$ws.Range("D12").Value = '5.21'
$ws.Range("E12").Value = '  +2.67%  '
$ws.Range("D13").Value = '0.347'
$ws.Range("E13").Value = '  +2.92%  '
$ws.Range("D14").Value = '26.40'
$ws.Range("E14").Value = '  +1.12%  '
$ws.Range("E15").Value = '  +1.57%  '
$ws.Range("D16").Value = '2.845.79'
$ws.Range("E16").Value = '  -0.32%  '
$ws.Range("D17").Value = '61.234.04'
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("D18").Value = '2.411.00'
$ws.Range("E18").Value = '  +0.38%  '
$ws.Range("D19").Value = '8.05'
$ws.Range("E19").Value = '  +0.32%  '
$ws.Range("D20").Value = '10.64'
$ws.Range("E20").Value = '  +0.81%  '
$ws.Range("D21").Value = '323.61'
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("E22").Value = '  +1.00%  '
$ws.Range("D23").Value = '6.09'
$ws.Range("E23").Value = '  -1.64%  '
$ws.Range("D24").Value = '1.97'
$ws.Range("E24").Value = '  +6.78%  '
$ws.Range("E25").Value = '  -0.29%  '
$ws.Range("D26").Value = '65.02'
$ws.Range("E26").Value = '  +0.90%  '
$ws.Range("D27").Value = '611.55'
$ws.Range("E27").Value = '  +4.59%  '
$ws.Range("E28").Value = '  -0.21%  '
$ws.Range("E29").Value = '  +2.01%  '
$ws.Range("D30").Value = '2.521.01'
$ws.Range("E30").Value = '  -0.86%  '
$ws.Range("D31").Value = '8.01'
$ws.Range("E31").Value = '  +0.43%  '
$ws.Range("D32").Value = '1.37'
$ws.Range("E32").Value = '  +2.61%  '
$ws.Range("D33").Value = '1.80'
$ws.Range("E33").Value = '  -0.83%  '
$ws.Range("E34").Value = '  +0.36%  '
$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D35").Value = '0.997'
$ws.Range("E35").Value = '  -0.66%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '1.46'
$ws.Range("E36").Value = '  +3.14%  '
$ws.Range("D37").Value = '152.72'
$ws.Range("E37").Value = '  +0.79%  '
$ws.Range("D38").Value = '0.372'
$ws.Range("E38").Value = '  +1.00%  '
$ws.Range("D39").Value = '4.63'
$ws.Range("E39").Value = '  +0.89%  '
$ws.Range("D40").Value = '5.31'
$ws.Range("E40").Value = '  +3.35%  '
$ws.Range("D41").Value = '18.34'
$ws.Range("E41").Value = '  +0.58%  '
$ws.Range("E42").Value = '  +6.28%  '
$ws.Range("E43").Value = '  +2.58%  '
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("D45").Value = '41.96'
$ws.Range("E45").Value = '  +1.88%  '
$ws.Range("E46").Value = '  -3.59%  '
$ws.Range("D47").Value = '142.22'
$ws.Range("E47").Value = '  -0.60%  '
$ws.Range("E48").Value = '  +0.34%  '
$ws.Range("D49").Value = '19.89'
$ws.Range("E49").Value = '  +2.87%  '
$ws.Range("E50").Value = '  +1.14%  '
$ws.Range("D51").Value = '0.0509'
$ws.Range("E51").Value = '  +1.60%  '
